$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds text that looks like a date ("2024-07-16"); prefix with an
# apostrophe (quote-prefix) so Excel keeps it as literal text instead of
# auto-converting it to a date serial number, matching the source data.

# Row 2
$ws.Range("A2").Value = "'2024-07-16"
$ws.Range("C2").Value = 88206
$ws.Range("D2").Value = "BR0026113"
$ws.Range("E2").Value = "MANGUINHOS ADM DE BENS E CONSULTORI"
$ws.Range("F2").Value = "RJ"
$ws.Range("G2").Value = 280.26
$ws.Range("H2").Value = 70.06999999999999
$ws.Range("I2").Value = 350.33

# Row 3
$ws.Range("A3").Value = "'2024-07-16"
$ws.Range("C3").Value = 88208
$ws.Range("D3").Value = "BR0010977"
$ws.Range("E3").Value = "AMB EMPREENDIMENTOS IMOBIL. LTDA"
$ws.Range("F3").Value = "SC"
$ws.Range("G3").Value = 2125.02
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 2125.02

# Row 4
$ws.Range("A4").Value = "'2024-07-16"
$ws.Range("C4").Value = 88211
$ws.Range("D4").Value = "BR0025869"
$ws.Range("E4").Value = "CONDOMINIO SOBERANE RESIDENCE, CORP"
$ws.Range("F4").Value = "AM"
$ws.Range("G4").Value = 793.61
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 793.61

# Row 5
$ws.Range("A5").Value = "'2024-07-16"
$ws.Range("C5").Value = 88212
$ws.Range("D5").Value = "BR0015419"
$ws.Range("E5").Value = "SPE SAUDE PRIMARIA BH S/A"
$ws.Range("F5").Value = "MG"
$ws.Range("G5").Value = 4397.29
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 4397.29
